$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.192.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.798.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4503"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +19.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07549"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.239"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.245"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.800.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.92%  "

$ws.Range("E17").Value = "  +3.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06679"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.414"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.203.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.418"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.004.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.289"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -13.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.931"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09450"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02380"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6736"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06278"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2169"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.190"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.483"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.217"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.156"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.865"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6114"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.035"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.168"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "
